$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the MuSCs-as-sender rows (rows 6-9); only the FAPs-as-sender
# rows (2-5) remain, now carrying updated TPM-derived values.
$ws.Rows("6:9").Delete()

# Row 2: FAPs -> Slitrk6/Ptprs -> ECs
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 3.556762333333333
$ws.Range("N2").Value = 10.670287
$ws.Range("O2").Value = 0.04280930450251701
$ws.Range("P2").Value = 0.04280930450251701
$ws.Range("Q2").Value = 6.571163463156223
$ws.Range("R2").Value = 59.140471168406
$ws.Range("S2").Value = 0.04280930450251701
$ws.Range("T2").Value = 0.04280930450251701

# Row 3: FAPs -> Slitrk6/Ptprs -> FAPs
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.5686906263805706
$ws.Range("P3").Value = 0.5686906263805704
$ws.Range("Q3").Value = 87.29315062083556
$ws.Range("R3").Value = 785.6383555875201
$ws.Range("S3").Value = 0.5686906263805706
$ws.Range("T3").Value = 0.5686906263805704

# Row 4: FAPs -> Slitrk6/Ptprs -> MuSCs
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 24.53173066666666
$ws.Range("N4").Value = 73.595192
$ws.Range("O4").Value = 0.2952646900921413
$ws.Range("P4").Value = 0.2952646900921412
$ws.Range("Q4").Value = 45.32268314192178
$ws.Range("R4").Value = 407.904148277296
$ws.Range("S4").Value = 0.2952646900921413
$ws.Range("T4").Value = 0.2952646900921412

# Row 5: FAPs -> Slitrk6/Ptprs -> Resolving-Mac
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("M5").Value = 7.746355333333334
$ws.Range("N5").Value = 23.239066
$ws.Range("O5").Value = 0.09323537902477132
$ws.Range("P5").Value = 0.0932353790247713
$ws.Range("Q5").Value = 14.31148959883422
$ws.Range("R5").Value = 128.803406389508
$ws.Range("S5").Value = 0.09323537902477132
$ws.Range("T5").Value = 0.0932353790247713
